# Apply the changes described by the diff:
# - Fill in the "d2" sheet's previously-blank (0) numeric cells with their
#   real values.
# - Move the active/selected tab from "d1" to "d2".
# - Update each sheet's cell selection to match the new state.

$wb = $excel.ActiveWorkbook

# --- d1 sheet: it stops being the active tab; just move its selection ---
$ws1 = $wb.Worksheets.Item("d1")
$ws1.Activate() | Out-Null
$ws1.Range("E29").Select() | Out-Null

# --- d2 sheet: fill in values, then make it the active tab ---
$ws2 = $wb.Worksheets.Item("d2")
$ws2.Activate() | Out-Null

$ws2.Range("A3").Value = 786.30769230769226
$ws2.Range("B3").Value = 947

$ws2.Range("A5").Value = 955.19230769230774
$ws2.Range("B5").Value = 1702.1538461538462

$ws2.Range("A7").Value = 1097.794117647059
$ws2.Range("B7").Value = 1266.610859728507

$ws2.Range("A34").Value = 3730.7692307692309
$ws2.Range("B34").Value = 30293.538461538461
$ws2.Range("C34").Value = 319.99450549450546

$ws2.Range("A36:B36").ClearFormats() | Out-Null
$ws2.Range("A36").Value = 11.076923076923077
$ws2.Range("B36").Value = 17.923076923076923
$ws2.Range("C36").Value = 1.9285714285714286

$ws2.Range("H20").Select() | Out-Null
